$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RYAAY")

  $ws.Cells.Item(8, 4).Value = 8023400
  $ws.Cells.Item(8, 5).Value = 7458800
  $ws.Cells.Item(8, 6).Value = 7333100
  $ws.Cells.Item(8, 7).Value = 6343700
  $ws.Cells.Item(8, 8).Value = 5651100
  $ws.Cells.Item(8, 9).Value = 5479800
  $ws.Cells.Item(8, 10).Value = 4925800
  $ws.Cells.Item(9, 4).Value = 4863600
  $ws.Cells.Item(9, 5).Value = 4665800
  $ws.Cells.Item(9, 6).Value = 4709700
  $ws.Cells.Item(9, 7).Value = 4346900
  $ws.Cells.Item(9, 8).Value = 4175700
  $ws.Cells.Item(9, 9).Value = 4451900
  $ws.Cells.Item(9, 10).Value = 3957300
  $ws.Cells.Item(10, 4).Value = 3159700
  $ws.Cells.Item(10, 5).Value = 2793000
  $ws.Cells.Item(10, 6).Value = 2623400
  $ws.Cells.Item(10, 7).Value = 1996800
  $ws.Cells.Item(10, 8).Value = 1475400
  $ws.Cells.Item(10, 9).Value = 1027900
  $ws.Cells.Item(10, 10).Value = 968500
  $ws.Cells.Item(17, 4).Value = 6152700
  $ws.Cells.Item(17, 5).Value = 5737600
  $ws.Cells.Item(17, 6).Value = 5694900
  $ws.Cells.Item(17, 7).Value = 5173600
  $ws.Cells.Item(17, 8).Value = 4912200
  $ws.Cells.Item(17, 9).Value = 4674000
  $ws.Cells.Item(17, 10).Value = 4159200
  $ws.Cells.Item(18, 4).Value = 1870700
  $ws.Cells.Item(18, 5).Value = 1721100
  $ws.Cells.Item(18, 6).Value = 1638200
  $ws.Cells.Item(18, 7).Value = 1170100
  $ws.Cells.Item(18, 8).Value = 738900
  $ws.Cells.Item(18, 9).Value = 805800
  $ws.Cells.Item(18, 10).Value = 766500
  $ws.Cells.Item(20, 4).Value = 4600
  $ws.Cells.Item(20, 5).Value = 3900
  $ws.Cells.Item(20, 6).Value = 373300
  $ws.Cells.Item(20, 7).Value = 15300
  $ws.Cells.Item(20, 8).Value = 16900
  $ws.Cells.Item(20, 9).Value = 35900
  $ws.Cells.Item(20, 10).Value = 66300
  $ws.Cells.Item(21, 4).Value = 2506200
  $ws.Cells.Item(21, 5).Value = 2284600
  $ws.Cells.Item(21, 6).Value = 2492100
  $ws.Cells.Item(21, 7).Value = 1610200
  $ws.Cells.Item(21, 8).Value = 1151600
  $ws.Cells.Item(21, 9).Value = 1212400
  $ws.Cells.Item(21, 10).Value = 1180600
  $ws.Cells.Item(22, 4).Value = 67400
  $ws.Cells.Item(22, 5).Value = 75400
  $ws.Cells.Item(22, 6).Value = 79500
  $ws.Cells.Item(22, 7).Value = 83100
  $ws.Cells.Item(22, 8).Value = 92300
  $ws.Cells.Item(22, 9).Value = 111400
  $ws.Cells.Item(22, 10).Value = 122600
  $ws.Cells.Item(23, 4).Value = 1807900
  $ws.Cells.Item(23, 5).Value = 1649700
  $ws.Cells.Item(23, 6).Value = 1932000
  $ws.Cells.Item(23, 7).Value = 1102200
  $ws.Cells.Item(23, 8).Value = 663500
  $ws.Cells.Item(23, 9).Value = 730300
  $ws.Cells.Item(23, 10).Value = 710200
  $ws.Cells.Item(24, 4).Value = 180800
  $ws.Cells.Item(24, 5).Value = 173200
  $ws.Cells.Item(24, 6).Value = 182700
  $ws.Cells.Item(24, 7).Value = 129800
  $ws.Cells.Item(24, 8).Value = 77000
  $ws.Cells.Item(24, 9).Value = 91600
  $ws.Cells.Item(24, 10).Value = 81500
  $ws.Cells.Item(26, 4).Value = 1627100
  $ws.Cells.Item(26, 5).Value = 1476400
  $ws.Cells.Item(26, 6).Value = 1749300
  $ws.Cells.Item(26, 7).Value = 972400
  $ws.Cells.Item(26, 8).Value = 586600
  $ws.Cells.Item(26, 9).Value = 638700
  $ws.Cells.Item(26, 10).Value = 628800
  $ws.Cells.Item(27, 4).Value = 1627100
  $ws.Cells.Item(27, 5).Value = 1476400
  $ws.Cells.Item(27, 6).Value = 1749300
  $ws.Cells.Item(27, 7).Value = 972400
  $ws.Cells.Item(27, 8).Value = 586600
  $ws.Cells.Item(27, 9).Value = 638700
  $ws.Cells.Item(27, 10).Value = 628800
  $ws.Cells.Item(32, 4).Value = -4600
  $ws.Cells.Item(32, 5).Value = -3900
  $ws.Cells.Item(32, 6).Value = -373300
  $ws.Cells.Item(32, 7).Value = -15300
  $ws.Cells.Item(32, 8).Value = -16900
  $ws.Cells.Item(32, 9).Value = -35900
  $ws.Cells.Item(32, 10).Value = -66300
  $ws.Cells.Item(33, 4).Value = 1627100
  $ws.Cells.Item(33, 5).Value = 1476400
  $ws.Cells.Item(33, 6).Value = 1749300
  $ws.Cells.Item(33, 7).Value = 972400
  $ws.Cells.Item(33, 8).Value = 586600
  $ws.Cells.Item(33, 9).Value = 638700
  $ws.Cells.Item(33, 10).Value = 628800
  $ws.Cells.Item(35, 4).Value = 1627100
  $ws.Cells.Item(35, 5).Value = 1476400
  $ws.Cells.Item(35, 6).Value = 1749300
  $ws.Cells.Item(35, 7).Value = 972400
  $ws.Cells.Item(35, 8).Value = 586600
  $ws.Cells.Item(35, 9).Value = 638700
  $ws.Cells.Item(35, 10).Value = 628800
  $ws.Cells.Item(41, 4).Value = 1699800
  $ws.Cells.Item(41, 5).Value = 1373300
  $ws.Cells.Item(41, 6).Value = 1412800
  $ws.Cells.Item(41, 7).Value = 1329100
  $ws.Cells.Item(41, 8).Value = 1941200
  $ws.Cells.Item(41, 9).Value = 1392300
  $ws.Cells.Item(41, 10).Value = 3038700
  $ws.Cells.Item(42, 4).Value = 2390400
  $ws.Cells.Item(42, 5).Value = 3258800
  $ws.Cells.Item(42, 6).Value = 3435900
  $ws.Cells.Item(42, 7).Value = 4044300
  $ws.Cells.Item(42, 8).Value = 1681100
  $ws.Cells.Item(42, 9).Value = 2573200
  $ws.Cells.Item(42, 10).Value = 866400
  $ws.Cells.Item(43, 4).Value = 65000
  $ws.Cells.Item(43, 5).Value = 62000
  $ws.Cells.Item(43, 6).Value = 78000
  $ws.Cells.Item(43, 7).Value = 73700
  $ws.Cells.Item(43, 8).Value = 69300
  $ws.Cells.Item(43, 9).Value = 66100
  $ws.Cells.Item(43, 10).Value = 63300
  $ws.Cells.Item(44, 4).Value = 4200
  $ws.Cells.Item(44, 5).Value = 3500
  $ws.Cells.Item(44, 6).Value = 3700
  $ws.Cells.Item(44, 7).Value = 2400
  $ws.Cells.Item(44, 8).Value = 2800
  $ws.Cells.Item(44, 9).Value = 3000
  $ws.Cells.Item(44, 10).Value = 3100
  $ws.Cells.Item(45, 4).Value = 540700
  $ws.Cells.Item(45, 5).Value = 582500
  $ws.Cells.Item(45, 6).Value = 479300
  $ws.Cells.Item(45, 7).Value = 993000
  $ws.Cells.Item(45, 8).Value = 170100
  $ws.Cells.Item(45, 9).Value = 188200
  $ws.Cells.Item(45, 10).Value = 377300
  $ws.Cells.Item(46, 4).Value = 4700000
  $ws.Cells.Item(46, 5).Value = 5280200
  $ws.Cells.Item(46, 6).Value = 5409700
  $ws.Cells.Item(46, 7).Value = 6442500
  $ws.Cells.Item(46, 8).Value = 3864500
  $ws.Cells.Item(46, 9).Value = 4222700
  $ws.Cells.Item(46, 10).Value = 4348800
  $ws.Cells.Item(47, 7).Value = 416300
  $ws.Cells.Item(47, 8).Value = 292100
  $ws.Cells.Item(47, 9).Value = 248200
  $ws.Cells.Item(47, 10).Value = 168000
  $ws.Cells.Item(48, 4).Value = 9114400
  $ws.Cells.Item(48, 5).Value = 8093800
  $ws.Cells.Item(48, 6).Value = 7025300
  $ws.Cells.Item(48, 7).Value = 6138500
  $ws.Cells.Item(48, 8).Value = 5677600
  $ws.Cells.Item(48, 9).Value = 5504800
  $ws.Cells.Item(48, 10).Value = 5526000
  $ws.Cells.Item(49, 4).Value = 52500
  $ws.Cells.Item(49, 5).Value = 52500
  $ws.Cells.Item(49, 6).Value = 52500
  $ws.Cells.Item(49, 7).Value = 52500
  $ws.Cells.Item(49, 8).Value = 52500
  $ws.Cells.Item(49, 9).Value = 52500
  $ws.Cells.Item(49, 10).Value = 52500
  $ws.Cells.Item(52, 4).Value = 2900
  $ws.Cells.Item(52, 5).Value = 25800
  $ws.Cells.Item(52, 6).Value = 99300
  $ws.Cells.Item(52, 7).Value = 622100
  $ws.Cells.Item(52, 8).Value = 400
  $ws.Cells.Item(52, 9).Value = 5700
  $ws.Cells.Item(52, 10).Value = 3700
  $ws.Cells.Item(54, 4).Value = 13869800
  $ws.Cells.Item(54, 5).Value = 13452300
  $ws.Cells.Item(54, 6).Value = 12586800
  $ws.Cells.Item(54, 7).Value = 13671900
  $ws.Cells.Item(54, 8).Value = 9887100
  $ws.Cells.Item(54, 9).Value = 10034000
  $ws.Cells.Item(54, 10).Value = 10099000
  $ws.Cells.Item(57, 4).Value = 280000
  $ws.Cells.Item(57, 5).Value = 330000
  $ws.Cells.Item(57, 6).Value = 258700
  $ws.Cells.Item(57, 7).Value = 220500
  $ws.Cells.Item(57, 8).Value = 168300
  $ws.Cells.Item(57, 9).Value = 155200
  $ws.Cells.Item(57, 10).Value = 203300
  $ws.Cells.Item(58, 4).Value = 487600
  $ws.Cells.Item(58, 5).Value = 511500
  $ws.Cells.Item(58, 6).Value = 504800
  $ws.Cells.Item(58, 7).Value = 448300
  $ws.Cells.Item(58, 8).Value = 525000
  $ws.Cells.Item(58, 9).Value = 448700
  $ws.Cells.Item(58, 10).Value = 413300
  $ws.Cells.Item(59, 4).Value = 3061600
  $ws.Cells.Item(59, 5).Value = 2537700
  $ws.Cells.Item(59, 6).Value = 3017000
  $ws.Cells.Item(59, 7).Value = 3085400
  $ws.Cells.Item(59, 8).Value = 1858700
  $ws.Cells.Item(59, 9).Value = 1541100
  $ws.Cells.Item(59, 10).Value = 1419800
  $ws.Cells.Item(60, 4).Value = 3829200
  $ws.Cells.Item(60, 5).Value = 3379200
  $ws.Cells.Item(60, 6).Value = 3780500
  $ws.Cells.Item(60, 7).Value = 3754200
  $ws.Cells.Item(60, 8).Value = 2552000
  $ws.Cells.Item(60, 9).Value = 2144900
  $ws.Cells.Item(60, 10).Value = 2036400
  $ws.Cells.Item(61, 4).Value = 3958800
  $ws.Cells.Item(61, 5).Value = 4407800
  $ws.Cells.Item(61, 6).Value = 4009000
  $ws.Cells.Item(61, 7).Value = 4523900
  $ws.Cells.Item(61, 8).Value = 2934800
  $ws.Cells.Item(61, 9).Value = 3476400
  $ws.Cells.Item(61, 10).Value = 3654100
  $ws.Cells.Item(62, 4).Value = 1067700
  $ws.Cells.Item(62, 5).Value = 702700
  $ws.Cells.Item(62, 6).Value = 761700
  $ws.Cells.Item(62, 7).Value = 866500
  $ws.Cells.Item(62, 8).Value = 713700
  $ws.Cells.Item(62, 9).Value = 740800
  $ws.Cells.Item(62, 10).Value = 698400
  $ws.Cells.Item(66, 4).Value = 8855800
  $ws.Cells.Item(66, 5).Value = 8489800
  $ws.Cells.Item(66, 6).Value = 8551200
  $ws.Cells.Item(66, 7).Value = 9144600
  $ws.Cells.Item(66, 8).Value = 6200500
  $ws.Cells.Item(66, 9).Value = 6362100
  $ws.Cells.Item(66, 10).Value = 6388900
  $ws.Cells.Item(72, 4).Value = 4602600
  $ws.Cells.Item(72, 5).Value = 3898200
  $ws.Cells.Item(72, 6).Value = 3565200
  $ws.Cells.Item(72, 7).Value = 3368800
  $ws.Cells.Item(72, 8).Value = 2980000
  $ws.Cells.Item(72, 9).Value = 2889200
  $ws.Cells.Item(72, 10).Value = 2796400
  $ws.Cells.Item(76, 4).Value = 5014100
  $ws.Cells.Item(76, 5).Value = 4962600
  $ws.Cells.Item(76, 6).Value = 4035600
  $ws.Cells.Item(76, 7).Value = 4527300
  $ws.Cells.Item(76, 8).Value = 3686600
  $ws.Cells.Item(76, 9).Value = 3671800
  $ws.Cells.Item(76, 10).Value = 3710100
  $ws.Cells.Item(81, 4).Value = 1627100
  $ws.Cells.Item(81, 5).Value = 1476400
  $ws.Cells.Item(81, 6).Value = 1749300
  $ws.Cells.Item(81, 7).Value = 972400
  $ws.Cells.Item(81, 8).Value = 586600
  $ws.Cells.Item(81, 9).Value = 638700
  $ws.Cells.Item(81, 10).Value = 628800
  $ws.Cells.Item(83, 4).Value = 629400
  $ws.Cells.Item(83, 5).Value = 558200
  $ws.Cells.Item(83, 6).Value = 479400
  $ws.Cells.Item(83, 7).Value = 423800
  $ws.Cells.Item(83, 8).Value = 394700
  $ws.Cells.Item(83, 9).Value = 369800
  $ws.Cells.Item(83, 10).Value = 346900
  $ws.Cells.Item(89, 4).Value = 2505600
  $ws.Cells.Item(89, 5).Value = 2162300
  $ws.Cells.Item(89, 6).Value = 2071500
  $ws.Cells.Item(89, 7).Value = 1895500
  $ws.Cells.Item(89, 8).Value = 1172000
  $ws.Cells.Item(89, 9).Value = 1148400
  $ws.Cells.Item(89, 10).Value = 1144800
  $ws.Cells.Item(91, 4).Value = -1650000
  $ws.Cells.Item(91, 5).Value = -1626700
  $ws.Cells.Item(91, 6).Value = -1366200
  $ws.Cells.Item(91, 7).Value = -884700
  $ws.Cells.Item(91, 8).Value = -567500
  $ws.Cells.Item(91, 9).Value = -348600
  $ws.Cells.Item(91, 10).Value = -356300
  $ws.Cells.Item(94, 4).Value = -807200
  $ws.Cells.Item(94, 5).Value = -1448300
  $ws.Cells.Item(94, 6).Value = -318200
  $ws.Cells.Item(94, 7).Value = -3240500
  $ws.Cells.Item(94, 8).Value = 337400
  $ws.Cells.Item(94, 9).Value = -2043700
  $ws.Cells.Item(94, 10).Value = -208000
  $ws.Cells.Item(96, 7).Value = -583800
  $ws.Cells.Item(96, 9).Value = -551500
  $ws.Cells.Item(100, 4).Value = -1372000
  $ws.Cells.Item(100, 5).Value = -753500
  $ws.Cells.Item(100, 6).Value = -1669600
  $ws.Cells.Item(100, 7).Value = 733000
  $ws.Cells.Item(100, 8).Value = -960500
  $ws.Cells.Item(100, 9).Value = -751100
  $ws.Cells.Item(100, 10).Value = -173800
  $ws.Cells.Item(102, 4).Value = 326500
  $ws.Cells.Item(102, 5).Value = -39500
  $ws.Cells.Item(102, 6).Value = 83700
  $ws.Cells.Item(102, 7).Value = -612000
  $ws.Cells.Item(102, 8).Value = 548900
  $ws.Cells.Item(102, 9).Value = -1646400
  $ws.Cells.Item(102, 10).Value = 763000
